# "new data morning 0708/2024"
# Replace the old placeholder rows (A2:C17, mostly empty/blank-styled cells)
# with a fresh list of 20 Facebook post URLs in column A, keeping the
# "link" header in A1.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Drop the column outline grouping (was outlineLevelCol="2" on sheetFormatPr).
$ws.Columns("A:A").OutlineLevel = 0

# Remove the old data rows entirely (2:17) — this also removes the stray
# B/C column cells that lived on rows 15-17 and the per-cell styling (s="2"/
# s="3"/s="4") those rows used to carry.
$ws.Rows("2:17").Delete()

$ws.Range("A1").Value = "link"

$urls = @(
    "https://www.facebook.com/viettan/posts/926081452895196?ref=embed_post",
    "https://www.facebook.com/chantroimoimedia/posts/903821068439059?ref=embed_post",
    "https://www.facebook.com/chantroimoimedia/posts/903153778505788?ref=embed_post",
    "https://www.facebook.com/chantroimoimedia/posts/902629631891536?ref=embed_post",
    "https://www.facebook.com/chantroimoimedia/posts/902511645236668?ref=embed_post",
    "https://www.facebook.com/chantroimoimedia/posts/902438565243976?ref=embed_post",
    "https://www.facebook.com/chantroimoimedia/posts/902270611927438?ref=embed_post",
    "https://www.facebook.com/chantroimoimedia/posts/902213735266459?ref=embed_post",
    "https://www.facebook.com/chantroimoimedia/posts/902185381935961?ref=embed_post",
    "https://www.facebook.com/viettan/posts/926752122828129?ref=embed_post",
    "https://www.facebook.com/viettan/posts/926457522857589?ref=embed_post",
    "https://www.facebook.com/viettan/posts/926382092865132?ref=embed_post",
    "https://www.facebook.com/viettan/posts/926308382872503?ref=embed_post",
    "https://www.facebook.com/viettan/posts/926251742878167?ref=embed_post",
    "https://www.facebook.com/viettan/posts/926178166218858?ref=embed_post",
    "https://www.facebook.com/viettan/posts/926081452895196?ref=embed_post",
    "https://www.facebook.com/viettan/posts/926045759565432?ref=embed_post",
    "https://www.facebook.com/viettan/posts/925408646295810?ref=embed_post",
    "https://www.facebook.com/viettan/posts/925780072925334?ref=embed_post",
    "https://www.facebook.com/viettan/posts/926976216139053?ref=embed_post"
)

for ($i = 0; $i -lt $urls.Length; $i++) {
    $ws.Cells.Item($i + 2, 1).Value = $urls[$i]
}

# Matches the saved selection left behind in the authored workbook.
$ws.Range("C14").Select() | Out-Null

